# Weekly fruit/vegetable price feed update.
# Two new weekly price records are prepended into the "Cebollin baby"
# data block (rows 43-63), pushing the existing rows down:
#   - a brand-new row is inserted at row 43
#   - a second brand-new row is inserted at row 48 (after the shift)
# All previously-existing rows 43-63 end up shifted down (by 1 for the
# first four rows, by 2 for the rest), ending at rows 44-47 and 49-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new row at 43 - shifts old row 43 (and below) down to 44.
$ws.Rows(43).Insert()

# Insert second new row at 48 - shifts old row 47 (now at row 47) and
# below down to 49+.
$ws.Rows(48).Insert()

# --- Fill the brand-new row 43 ---
$ws.Cells.Item(43,1).Value2  = 1
$ws.Cells.Item(43,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(43,4).Value2  = 44574
$ws.Cells.Item(43,5).Value2  = 15
$ws.Cells.Item(43,6).Value2  = 100112038
$ws.Cells.Item(43,7).Value2  = "Cebollín baby"
$ws.Cells.Item(43,8).Value2  = "Sin especificar"
$ws.Cells.Item(43,9).Value2  = "Primera"
$ws.Cells.Item(43,10).Value2 = 200
$ws.Cells.Item(43,11).Value2 = 5000
$ws.Cells.Item(43,12).Value2 = 5500
$ws.Cells.Item(43,13).Value2 = 5250
$ws.Cells.Item(43,14).Value2 = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(43,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(43,16).Value2 = 2625
$ws.Cells.Item(43,17).Value2 = 2
$ws.Cells.Item(43,18).Value2 = "Hortaliza"

# --- Fill the brand-new row 48 ---
$ws.Cells.Item(48,1).Value2  = 1
$ws.Cells.Item(48,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(48,4).Value2  = 44230
$ws.Cells.Item(48,5).Value2  = 15
$ws.Cells.Item(48,6).Value2  = 100112038
$ws.Cells.Item(48,7).Value2  = "Cebollín baby"
$ws.Cells.Item(48,8).Value2  = "Sin especificar"
$ws.Cells.Item(48,9).Value2  = "Primera"
$ws.Cells.Item(48,10).Value2 = 250
$ws.Cells.Item(48,11).Value2 = 5500
$ws.Cells.Item(48,12).Value2 = 6000
$ws.Cells.Item(48,13).Value2 = 5750
$ws.Cells.Item(48,14).Value2 = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(48,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(48,16).Value2 = 2875
$ws.Cells.Item(48,17).Value2 = 2
$ws.Cells.Item(48,18).Value2 = "Hortaliza"
